# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" sheets, which contain the same underlying data.

$wb = $excel.ActiveWorkbook

# Row -> new F value (column F = "想去人数")
$updates = @{
    7  = 1282
    10 = 404
    12 = 161
    17 = 309
    19 = 1743
    21 = 107
    26 = 4210
    29 = 1097
    32 = 569
    34 = 271
    36 = 142
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
